$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4962.2
$ws.Range("I62").Value = 3148.75
$ws.Range("J62").Value = 6171.1665
$ws.Range("K62").Value = 3148.75
$ws.Range("L62").Value = 6171.1665
$ws.Range("M62").Value = -2524.75
$ws.Range("N62").Value = -7419.1665
$ws.Range("H65").Value = 4962.2
$ws.Range("I65").Value = 3148.75
$ws.Range("J65").Value = 6171.1665
$ws.Range("K65").Value = 15743.75
$ws.Range("L65").Value = 30855.8325
$ws.Range("M65").Value = -12623.75
$ws.Range("N65").Value = -37095.8325
$ws.Range("H76").Value = 3051.84
$ws.Range("I76").Value = 2994
$ws.Range("J76").Value = 3200.5715
$ws.Range("K76").Value = 2994
$ws.Range("L76").Value = 3200.5715
$ws.Range("M76").Value = -2679
$ws.Range("N76").Value = -3830.5715
$ws.Range("H79").Value = 3051.84
$ws.Range("I79").Value = 2994
$ws.Range("J79").Value = 3200.5715
$ws.Range("K79").Value = 2994
$ws.Range("L79").Value = 3200.5715
$ws.Range("M79").Value = -1902
$ws.Range("N79").Value = -5384.5715
$ws.Range("H98").Value = 1943.3572
$ws.Range("I98").Value = 1387
$ws.Range("J98").Value = 2499.7144
$ws.Range("K98").Value = 1387
$ws.Range("L98").Value = 2499.7144
$ws.Range("M98").Value = 111
$ws.Range("N98").Value = -5495.7144
$ws.Range("H122").Value = 1943.3572
$ws.Range("I122").Value = 1387
$ws.Range("J122").Value = 2499.7144
$ws.Range("K122").Value = 4161
$ws.Range("L122").Value = 7499.1432
$ws.Range("M122").Value = -1711
$ws.Range("N122").Value = -12399.1432
$ws.Range("H137").Value = 1352
$ws.Range("I137").Value = 893.13635
$ws.Range("J137").Value = 3875.75
$ws.Range("K137").Value = 2679.40905
$ws.Range("L137").Value = 11627.25
$ws.Range("M137").Value = -129.4090500000002
$ws.Range("N137").Value = -16727.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2840.7778
$ws.Range("I63").Value = 1927.8334
$ws.Range("J63").Value = 4666.6665
$ws.Range("K63").Value = 1927.8334
$ws.Range("L63").Value = 4666.6665
$ws.Range("M63").Value = -1241.8334
$ws.Range("N63").Value = -6038.6665
$ws.Range("H66").Value = 2840.7778
$ws.Range("I66").Value = 1927.8334
$ws.Range("J66").Value = 4666.6665
$ws.Range("K66").Value = 9639.166999999999
$ws.Range("L66").Value = 23333.3325
$ws.Range("M66").Value = -6207.166999999999
$ws.Range("N66").Value = -30197.3325
$ws.Range("H132").Value = 2390.5405
$ws.Range("I132").Value = 2213.25
$ws.Range("J132").Value = 2942.111
$ws.Range("K132").Value = 6639.75
$ws.Range("L132").Value = 8826.332999999999
$ws.Range("M132").Value = -4109.75
$ws.Range("N132").Value = -13886.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5956.852
$ws.Range("I134").Value = 5053.737
$ws.Range("J134").Value = 8101.75
$ws.Range("K134").Value = 15161.211
$ws.Range("L134").Value = 24305.25
$ws.Range("M134").Value = -12626.211
$ws.Range("N134").Value = -29375.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 833.3
$ws.Range("I16").Value = 747.875
$ws.Range("J16").Value = 1175
$ws.Range("K16").Value = 747.875
$ws.Range("L16").Value = 1175
$ws.Range("M16").Value = -460.875
$ws.Range("N16").Value = -1749
$ws.Range("H31").Value = 21277796
$ws.Range("I31").Value = 34483536
$ws.Range("J31").Value = 1884.6666
$ws.Range("K31").Value = 34483536
$ws.Range("L31").Value = 1884.6666
$ws.Range("M31").Value = -34483241
$ws.Range("N31").Value = -2474.6666
$ws.Range("H34").Value = 21277796
$ws.Range("I34").Value = 34483536
$ws.Range("J34").Value = 1884.6666
$ws.Range("K34").Value = 34483536
$ws.Range("L34").Value = 1884.6666
$ws.Range("M34").Value = -34483334
$ws.Range("N34").Value = -2288.6666
$ws.Range("H58").Value = 1392.9474
$ws.Range("I58").Value = 1228.2667
$ws.Range("J58").Value = 2010.5
$ws.Range("K58").Value = 1228.2667
$ws.Range("L58").Value = 2010.5
$ws.Range("M58").Value = -1025.2667
$ws.Range("N58").Value = -2416.5
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H113").Value = 833.3
$ws.Range("I113").Value = 747.875
$ws.Range("J113").Value = 1175
$ws.Range("K113").Value = 747.875
$ws.Range("L113").Value = 1175
$ws.Range("M113").Value = 1422.125
$ws.Range("N113").Value = -5515
$ws.Range("H134").Value = 1824.8334
$ws.Range("I134").Value = 1926.5652
$ws.Range("J134").Value = 1490.5714
$ws.Range("K134").Value = 5779.6956
$ws.Range("L134").Value = 4471.7142
$ws.Range("M134").Value = -3244.6956
$ws.Range("N134").Value = -9541.7142
$ws.Range("H135").Value = 34000
$ws.Range("J135").Value = 34000
$ws.Range("L135").Value = 34000
$ws.Range("N135").Value = -44140
$ws.Range("H136").Value = 1392.9474
$ws.Range("I136").Value = 1228.2667
$ws.Range("J136").Value = 2010.5
$ws.Range("K136").Value = 3684.800099999999
$ws.Range("L136").Value = 6031.5
$ws.Range("M136").Value = -1134.800099999999
$ws.Range("N136").Value = -11131.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 796.83
$ws.Range("I131").Value = 370.94446
$ws.Range("J131").Value = 890.3171
$ws.Range("K131").Value = 1112.83338
$ws.Range("L131").Value = 2670.9513
$ws.Range("M131").Value = 3927.16662
$ws.Range("N131").Value = -12750.9513

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4913.1313
$ws.Range("I70").Value = 4832
$ws.Range("J70").Value = 5069.154
$ws.Range("K70").Value = 4832
$ws.Range("L70").Value = 5069.154
$ws.Range("M70").Value = -4562
$ws.Range("N70").Value = -5609.154
$ws.Range("H73").Value = 4913.1313
$ws.Range("I73").Value = 4832
$ws.Range("J73").Value = 5069.154
$ws.Range("K73").Value = 4832
$ws.Range("L73").Value = 5069.154
$ws.Range("M73").Value = -3896
$ws.Range("N73").Value = -6941.154
$ws.Range("H126").Value = 2009.8
$ws.Range("I126").Value = 2012.25
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6036.75
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3566.75
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 3605.7368
$ws.Range("I132").Value = 2967.4
$ws.Range("J132").Value = 5999.5
$ws.Range("K132").Value = 8902.200000000001
$ws.Range("L132").Value = 17998.5
$ws.Range("M132").Value = -6372.200000000001
$ws.Range("N132").Value = -23058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2961.36
$ws.Range("I132").Value = 2342.5833
$ws.Range("J132").Value = 3532.5386
$ws.Range("K132").Value = 7027.749899999999
$ws.Range("L132").Value = 10597.6158
$ws.Range("M132").Value = -4497.749899999999
$ws.Range("N132").Value = -15657.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6456955
$ws.Range("I132").Value = 8700774
$ws.Range("K132").Value = 26102322
$ws.Range("M132").Value = -26099792
